# Applies the "Generate Report for Archive" update:
# Re-sorts rows 2-4 across all three sheets (Overview, zh-cn, de-de)
# so the 51eb981c record moves to row 2, 690dcbef to row 3, e133a15c to row 4.
# The 51eb981c record also flips status from "Ready for handoff" back to "In Translation".
# Row 5 (a9e883ab) is unchanged.

$wb = $excel.ActiveWorkbook

# --- Sheet "Overview" ---
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Hyperlinks.Delete()

# Row 2: 51eb981c-dc17-414f-bafa-f65f151ce654
$ws1.Range("B2").Value = "In Translation"
$ws1.Range("C2").Value = "In Translation"
$ws1.Range("D2").Value = "2016-03-25 07:53:58"
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a761c4671b6e439c8ffeba0b76565761882a850c/e2e/51eb981c-dc17-414f-bafa-f65f151ce654.md", "", "", "51eb981c-dc17-414f-bafa-f65f151ce654.md")

# Row 3: 690dcbef-f1d6-4296-ab6a-ed3e383eb003
$ws1.Range("B3").Value = "In Translation"
$ws1.Range("C3").Value = "In Translation"
$ws1.Range("D3").Value = "2016-03-25 07:52:19"
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/02accbff4377db5cbb4d463c8c8bc2770a9ce524/e2e/690dcbef-f1d6-4296-ab6a-ed3e383eb003.md", "", "", "690dcbef-f1d6-4296-ab6a-ed3e383eb003.md")

# Row 4: e133a15c-b1a5-41db-9e69-610819c60f4d
$ws1.Range("B4").Value = "In Translation"
$ws1.Range("C4").Value = "In Translation"
$ws1.Range("D4").Value = "2016-03-25 07:52:19"
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/02accbff4377db5cbb4d463c8c8bc2770a9ce524/e2e/e133a15c-b1a5-41db-9e69-610819c60f4d.md", "", "", "e133a15c-b1a5-41db-9e69-610819c60f4d.md")

# --- Sheet "zh-cn" ---
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Hyperlinks.Delete()

# Row 2: 51eb981c-dc17-414f-bafa-f65f151ce654
$ws2.Range("C2").Value = "In Translation"
$ws2.Range("E2").Value = "2016-03-25 07:53:49"
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a761c4671b6e439c8ffeba0b76565761882a850c/e2e/51eb981c-dc17-414f-bafa-f65f151ce654.md", "", "", "51eb981c-dc17-414f-bafa-f65f151ce654.md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e99eced8c0a9414d4ee16347533fdfb3355df856/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/51eb981c-dc17-414f-bafa-f65f151ce654.c6228fe864671987c0d0cc7531906fe94641eda8.zh-cn.xlf", "", "", "51eb981c-dc17-414f-bafa-f65f151ce654.c6228fe864671987c0d0cc7531906fe94641eda8.zh-cn.xlf")

# Row 3: 690dcbef-f1d6-4296-ab6a-ed3e383eb003
$ws2.Range("C3").Value = "In Translation"
$ws2.Range("E3").Value = "2016-03-25 07:52:06"
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/02accbff4377db5cbb4d463c8c8bc2770a9ce524/e2e/690dcbef-f1d6-4296-ab6a-ed3e383eb003.md", "", "", "690dcbef-f1d6-4296-ab6a-ed3e383eb003.md")
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0fa9243e62fde71fb39dbcf252cab93e4dc2f58e/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/690dcbef-f1d6-4296-ab6a-ed3e383eb003.31348ba87969569e75cd4b08fa060c15c4cdc9a1.zh-cn.xlf", "", "", "690dcbef-f1d6-4296-ab6a-ed3e383eb003.31348ba87969569e75cd4b08fa060c15c4cdc9a1.zh-cn.xlf")

# Row 4: e133a15c-b1a5-41db-9e69-610819c60f4d
$ws2.Range("C4").Value = "In Translation"
$ws2.Range("E4").Value = "2016-03-25 07:52:06"
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/02accbff4377db5cbb4d463c8c8bc2770a9ce524/e2e/e133a15c-b1a5-41db-9e69-610819c60f4d.md", "", "", "e133a15c-b1a5-41db-9e69-610819c60f4d.md")
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0fa9243e62fde71fb39dbcf252cab93e4dc2f58e/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/e133a15c-b1a5-41db-9e69-610819c60f4d.36be94687cdae3f623fe6203cf77d9203112376f.zh-cn.xlf", "", "", "e133a15c-b1a5-41db-9e69-610819c60f4d.36be94687cdae3f623fe6203cf77d9203112376f.zh-cn.xlf")

# --- Sheet "de-de" ---
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Hyperlinks.Delete()

# Row 2: 51eb981c-dc17-414f-bafa-f65f151ce654
$ws3.Range("C2").Value = "In Translation"
$ws3.Range("E2").Value = "2016-03-25 07:53:58"
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a761c4671b6e439c8ffeba0b76565761882a850c/e2e/51eb981c-dc17-414f-bafa-f65f151ce654.md", "", "", "51eb981c-dc17-414f-bafa-f65f151ce654.md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/30991a27a5173cb810002e05f8df9222416bd2ea/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/51eb981c-dc17-414f-bafa-f65f151ce654.c6228fe864671987c0d0cc7531906fe94641eda8.de-de.xlf", "", "", "51eb981c-dc17-414f-bafa-f65f151ce654.c6228fe864671987c0d0cc7531906fe94641eda8.de-de.xlf")

# Row 3: 690dcbef-f1d6-4296-ab6a-ed3e383eb003
$ws3.Range("C3").Value = "In Translation"
$ws3.Range("E3").Value = "2016-03-25 07:52:19"
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/02accbff4377db5cbb4d463c8c8bc2770a9ce524/e2e/690dcbef-f1d6-4296-ab6a-ed3e383eb003.md", "", "", "690dcbef-f1d6-4296-ab6a-ed3e383eb003.md")
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c1908ca9560353f434faac32bc5d57a98e77b2d2/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/690dcbef-f1d6-4296-ab6a-ed3e383eb003.31348ba87969569e75cd4b08fa060c15c4cdc9a1.de-de.xlf", "", "", "690dcbef-f1d6-4296-ab6a-ed3e383eb003.31348ba87969569e75cd4b08fa060c15c4cdc9a1.de-de.xlf")

# Row 4: e133a15c-b1a5-41db-9e69-610819c60f4d
$ws3.Range("C4").Value = "In Translation"
$ws3.Range("E4").Value = "2016-03-25 07:52:19"
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/02accbff4377db5cbb4d463c8c8bc2770a9ce524/e2e/e133a15c-b1a5-41db-9e69-610819c60f4d.md", "", "", "e133a15c-b1a5-41db-9e69-610819c60f4d.md")
$ws3.Hyperlinks.Add($ws3.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c1908ca9560353f434faac32bc5d57a98e77b2d2/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/e133a15c-b1a5-41db-9e69-610819c60f4d.36be94687cdae3f623fe6203cf77d9203112376f.de-de.xlf", "", "", "e133a15c-b1a5-41db-9e69-610819c60f4d.36be94687cdae3f623fe6203cf77d9203112376f.de-de.xlf")
